$d = $word.ActiveDocument

function Find-Range($startPos, $endPos, $text) {
    $r = $d.Range($startPos, $endPos)
    $found = $r.Find.Execute($text, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    return @{ found = $found; r = $r }
}

# ===================================================================
# Hunk 1: " 1" + "3" + "/" -> " 13/"  (merge 3 runs for the date " 13/05/2021")
# Runs "Thursday", "0", "5", "/2021" share identical formatting with
# " 1"+"3"+"/" so any touch coalesces the whole span; force-merge then
# re-split at the boundaries that must remain separate runs.
# ===================================================================
$whole0 = $d.Content
$f0 = $whole0.Find.Execute("Thursday 13/05/2021", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$locStart0 = $whole0.Start
$locEnd0 = $whole0.End

$r0 = Find-Range $locStart0 $locEnd0 " 13/"
$r0.r.Text = " XX/"
$r0.r.Text = " 13/"

$r0b = Find-Range $locStart0 $locEnd0 " 13/"
$r0b.r.Bold = 1
$r0b.r.Bold = 0

$r0c = Find-Range $locStart0 $locEnd0 "0"
$r0c.r.Bold = 1
$r0c.r.Bold = 0

$r0d = Find-Range $locStart0 $locEnd0 "5"
$r0d.r.Bold = 1
$r0d.r.Bold = 0

# ===================================================================
# Hunk 2: "Mia " + "Vasiliadis" (proofErr-wrapped) + " " -> "Mia Vasiliadis "
# ===================================================================
$d.Content.Find.Execute("Mia Vasiliadis ", $false, $false, $false, $false, $false, $true, 1, $false, "Mia Vasiliadis ", 2) | Out-Null

# ===================================================================
# Hunk 3: "Due date: Sunday of week " + "13(?):" ->
#         "Due date: " | "Friday" | " of week " | "13:"   (4 runs)
# ===================================================================
$whole = $d.Content
$f = $whole.Find.Execute("Due date: Sunday of week 13(?):", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$locStart = $whole.Start
$locEnd = $whole.End

# remove "(?)"
$res = Find-Range $locStart $locEnd "(?)"
$res.r.Text = ""
$locEnd = $locEnd - 3

# replace Sunday -> Friday, forcing a run split via a harmless bold toggle
$res2 = Find-Range $locStart $locEnd "Sunday"
$res2.r.Text = "Friday"
$res2.r.Bold = 1
$res2.r.Bold = 0

# force a run split before the trailing "13:"
$res3 = Find-Range $locStart $locEnd "13:"
$res3.r.Bold = 1
$res3.r.Bold = 0

# ===================================================================
# Hunk 4: "Sunday of week " + "13(?):" -> "Friday" | " of week " | "13:"  (3 runs)
# ===================================================================
$whole2 = $d.Content
$f2 = $whole2.Find.Execute("Final for A5 due: Sunday of week 13(?):", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$locStart2 = $whole2.Start
$locEnd2 = $whole2.End

$res4 = Find-Range $locStart2 $locEnd2 "(?)"
$res4.r.Text = ""
$locEnd2 = $locEnd2 - 3

$res5 = Find-Range $locStart2 $locEnd2 "Sunday"
$res5.r.Text = "Friday"
$res5.r.Bold = 1
$res5.r.Bold = 0

$res6 = Find-Range $locStart2 $locEnd2 "13:"
$res6.r.Bold = 1
$res6.r.Bold = 0

# ===================================================================
# Hunk 5: "Timeframe table" + " (Mia)" + "." -> "Timeframe table (Mia)."
# ===================================================================
$d.Content.Find.Execute("Timeframe table (Mia).", $false, $false, $false, $false, $false, $true, 1, $false, "Timeframe table (Mia).", 2) | Out-Null
